$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.815.22"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "'2.532.01"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'305.38"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'99.35"
$ws.Range("E6").Value = "  +4.74%  "
$ws.Range("D7").Value = "'0.583"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.547"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").Value = "'37.01"
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("D11").Value = "'0.0814"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "'7.74"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "'2.921.15"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "'2.564.39"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "'15.15"
$ws.Range("E16").Value = "  +6.53%  "
$ws.Range("D17").Value = "'0.870"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "'42.902.21"
$ws.Range("D19").Value = "'13.04"
$ws.Range("E19").Value = "  +2.94%  "
$ws.Range("D20").Value = "'0.0₃0983"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").Value = "'6.49"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").Value = "'71.55"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "'253.77"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").Value = "'2.93"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").Value = "'2.05"
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").Value = "'27.16"
$ws.Range("E26").Value = "  -5.48%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'10.46"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D29").Value = "'2.31"
$ws.Range("E29").Value = "  +8.35%  "
$ws.Range("D30").Value = "'38.69"
$ws.Range("E30").Value = "  +4.75%  "
$ws.Range("D31").Value = "'6.14"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").Value = "'158.57"
$ws.Range("E32").Value = "  +2.86%  "
$ws.Range("D33").Value = "'3.32"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").Value = "'2.11"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("D35").Value = "'0.0794"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").Value = "'2.63"
$ws.Range("E36").Value = "  -4.07%  "
$ws.Range("D37").Value = "'18.26"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").Value = "'0.115"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'24.30"
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.120"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").Value = "'3.45"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("B42").Value = "ApeXProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D42").Value = "'2.09"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'3.90"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "'0.0304"
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("D45").Value = "'0.998"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'2.042.03"
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").Value = "'85.87"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("D48").Value = "'9.00"
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("D49").Value = "'2.782.42"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.192"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'103.01"
$ws.Range("E51").Value = "  -3.87%  "
